# fix(data): robust leftover-window calculation per column
#
# The MSCI_ACWI series (column C) runs out of source data 39 months
# before the other series in this sheet. Those trailing rows were
# previously carrying stale/duplicated numbers; this edit removes them
# so each column's "leftover window" (trailing partial window of
# months) is computed independently from its own last valid period
# instead of assuming every column has data through the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-invalid MSCI_ACWI figures for the trailing months that
# no longer have real source data (rows 283-321 of column C).
$ws.Range("C283:C321").Clear() | Out-Null

# Leave the view where the analyst was looking when verifying the fix:
# scrolled down to, and focused on, the tail of the MSCI_ACWI column.
$win = $excel.ActiveWindow
$win.ScrollRow = 266
$win.ScrollColumn = 1
$ws.Range("F280").Select() | Out-Null
